$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (33) down onto the two
# new rows (34:35) first, so the new cells pick up style index 1 (the
# default centered data-row style) just like every other data row.
$ws.Range("A33:R33").Copy()
$ws.Range("A34:R35").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 34 - new local extreme: Лаишевский муниципальный район, female, 2020
$ws.Range("A34").Value = 92634000
$ws.Range("B34").Value = "Лаишевский муниципальный район"
$ws.Range("C34").Value = "female"
$ws.Range("D34").Value = 2020
$ws.Range("E34").Value = 0.06177
$ws.Range("F34").Value = 0.0672
$ws.Range("G34").Value = 0.0521
$ws.Range("H34").Value = 0.04675
$ws.Range("I34").Value = 0.04807
$ws.Range("J34").Value = 0.0706
$ws.Range("K34").Value = 0.0935
$ws.Range("L34").Value = 0.0807
$ws.Range("M34").Value = 0.0676
$ws.Range("N34").Value = 0.05908
$ws.Range("O34").Value = 0.1287
$ws.Range("P34").Value = 0.07776
$ws.Range("Q34").Value = 0.08466
$ws.Range("R34").Value = 0.06152

# Row 35 - same district, male, 2020
$ws.Range("A35").Value = 92634000
$ws.Range("B35").Value = "Лаишевский муниципальный район"
$ws.Range("C35").Value = "male"
$ws.Range("D35").Value = 2020
$ws.Range("E35").Value = 0.0675
$ws.Range("F35").Value = 0.07605
$ws.Range("G35").Value = 0.05786
$ws.Range("H35").Value = 0.04883
$ws.Range("I35").Value = 0.04376
$ws.Range("J35").Value = 0.0773
$ws.Range("K35").Value = 0.10565
$ws.Range("L35").Value = 0.08075
$ws.Range("M35").Value = 0.0722
$ws.Range("N35").Value = 0.05878
$ws.Range("O35").Value = 0.1194
$ws.Range("P35").Value = 0.07074
$ws.Range("Q35").Value = 0.07214
$ws.Range("R35").Value = 0.04916

# Match the author's final view state: scrolled down a bit further, with
# B38 (first empty cell below the new block) selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$null = $ws.Range("B38").Select()
